$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monsters")

# Row 3: fire_atonement (AQ) and water_atonement (AS) updated
$ws.Range("AQ3").Value = 0.18
$ws.Range("AS3").Value = 0.17

# Row 5: fire_atonement (AQ) updated
$ws.Range("AQ5").Value = 0.2

# Row 6: fire_atonement (AQ) updated
$ws.Range("AQ6").Value = 0.32

# Row 7: fire_atonement (AQ) updated
$ws.Range("AQ7").Value = 0.19

# Row 8: fire_atonement (AQ), ice_atonement (AR), water_atonement (AS) updated
$ws.Range("AQ8").Value = 0.45
$ws.Range("AR8").Value = 0.14
$ws.Range("AS8").Value = 0.16
